$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.046.87'
$ws.Range("E2").Value = '  -1.12%  '
$ws.Range("D3").Value = '1.916.86'
$ws.Range("E3").Value = '  -3.99%  '
$ws.Range("E4").Value = '  +0.13%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '239.06'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  -2.95%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '0.599'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  -5.08%  '
$ws.Range("E7").Value = '  +0.03%  '
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '54.97'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  -10.89%  '
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.357'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  -7.17%  '
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '55.07'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  -3.22%  '
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.0806'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  +3.95%  '
$ws.Range("E12").Value = '  -0.76%  '
$ws.Range("D13").Value = '2.203.91'
$ws.Range("E13").Value = '  -3.80%  '
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '0.796'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  -9.35%  '
# Rows 15 & 16: Avalanche/Chainlink swap positions with refreshed values
$cell = $ws.Range("B15")
$cell.Value = 'Chainlink'
$cell = $ws.Range("C15")
$cell.Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '12.99'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  -8.19%  '
$cell = $ws.Range("B16")
$cell.Value = 'Avalanche'
$cell = $ws.Range("C16")
$cell.Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '20.33'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  -12.19%  '

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '5.10'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  -7.15%  '
$ws.Range("D18").Value = '1.896.58'
$ws.Range("E18").Value = '  -5.76%  '
$ws.Range("D19").Value = '35.927.40'
$ws.Range("E19").Value = '  -1.15%  '
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '68.68'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  -4.17%  '
$ws.Range("D21").Value = '0.0₃0847'
$ws.Range("E21").Value = '  -2.63%  '
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '224.35'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  -4.28%  '
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '4.87'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  -7.74%  '
$ws.Range("E24").Value = '  +0.42%  '
$ws.Range("E25").Value = '  -4.22%  '
$ws.Range("E26").Value = '  -3.52%  '
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '9.13'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  -5.78%  '
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '161.97'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  +1.65%  '
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '18.90'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  -5.87%  '
$ws.Range("E30").Value = '  -17.59%  '
$ws.Range("E31").Value = '  -3.56%  '
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '1.10'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  -5.71%  '
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '4.55'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  -8.04%  '
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '0.0611'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  -0.57%  '
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '4.17'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  -6.10%  '
$ws.Range("E36").Value = '  +0.18%  '
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '1.80'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -1.76%  '
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '5.85'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  -10.25%  '
$ws.Range("E39").Value = '  -11.19%  '
$ws.Range("E40").Value = '  -13.14%  '
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '0.0949'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  -4.37%  '
$ws.Range("E42").Value = '  -1.48%  '
$ws.Range("E43").Value = '  -8.70%  '
$ws.Range("E44").Value = '  -4.27%  '
$ws.Range("D45").Value = '1.321.21'
$ws.Range("E45").Value = '  -2.41%  '
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '15.19'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  -8.78%  '
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  -10.00%  '
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '85.78'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  -7.06%  '
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '7.04'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  -6.97%  '
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '2.78'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  -2.71%  '
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '44.68'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  +0.55%  '
